# "slides for topic 2"
#
# Slide 1, shape id=3 ("Subtitle 2" / subTitle placeholder) is rewritten:
#   - paragraph "Andrew Beatty"            -> "Programming and Scripting"
#   - paragraph "Andrew.Beatty@gmit.ie"    -> split into two hyperlinked runs
#                                              "Andrew.Beatty@atu" + ".ie"
#                                              (hyperlink target unchanged)
#   - paragraph "Programming and Scripting" is dropped (its text moved to
#     paragraph 1 above, so the text frame goes from 3 paragraphs to 2)
#   - the body's normAutofit shrink (fontScale/lnSpcReduction) is cleared
#     back to a plain <a:normAutofit/> now that the text is shorter.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Subtitle 2")
$tf  = $shp.TextFrame
$tr  = $tf.TextRange

# Drop the autofit shrink recorded for the old (longer) text -- resetting
# AutoSize re-serialises normAutofit without fontScale/lnSpcReduction.
$tf.AutoSize = 2

# Remember the hyperlink target before the text underneath it is replaced.
$mailtoAddress = $tr.Paragraphs(2).Runs(1).ActionSettings.Item(1).Hyperlink.Address

# Rewrite the whole frame: "Programming and Scripting" becomes its own first
# paragraph, the (soon to be re-linked) email becomes the second -- this
# collapses the old 3-paragraph body down to 2 paragraphs in one shot.
$tr.Text = "Programming and Scripting" + [char]13 + "Andrew.Beatty@atu.ie"

# Re-apply the hyperlink to the new email paragraph, split across two runs
# exactly as authored: "Andrew.Beatty@atu" + ".ie".
$para2 = $tr.Paragraphs(2)
$run1  = $tr.Characters($para2.Start, 17)
$run2  = $tr.Characters($para2.Start + 17, 3)
$run1.ActionSettings.Item(1).Hyperlink.Address = $mailtoAddress
$run2.ActionSettings.Item(1).Hyperlink.Address = $mailtoAddress
